$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix ROI positions: shift x from 0.6 -> 0.7 and y from 0.6 -> 0.5 (keeping signs),
# applied to every row that references each of the four unique position strings.

$ws.Range("A2").Value = "(-0.7, 0.5)"
$ws.Range("A6").Value = "(-0.7, 0.5)"
$ws.Range("A10").Value = "(-0.7, 0.5)"
$ws.Range("A14").Value = "(-0.7, 0.5)"

$ws.Range("A3").Value = "(-0.7, -0.5)"
$ws.Range("A7").Value = "(-0.7, -0.5)"
$ws.Range("A11").Value = "(-0.7, -0.5)"
$ws.Range("A15").Value = "(-0.7, -0.5)"

$ws.Range("A4").Value = "(0.7, 0.5)"
$ws.Range("A8").Value = "(0.7, 0.5)"
$ws.Range("A12").Value = "(0.7, 0.5)"
$ws.Range("A16").Value = "(0.7, 0.5)"

$ws.Range("A5").Value = "(0.7, -0.5)"
$ws.Range("A9").Value = "(0.7, -0.5)"
$ws.Range("A13").Value = "(0.7, -0.5)"
$ws.Range("A17").Value = "(0.7, -0.5)"

# Update the active selection to match the last-saved cursor position (A14:A17).
$ws.Range("A14:A17").Select()
